$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A new observation (2026/02/27, 金, hour 4, rank 201) was recorded between the
# existing 2026/02/27 01:00 row and the 2026/12/29 block. Insert a new row at
# row 893 (pushing the old row 893 onward down to row 894, etc.) and fill it
# with the new data point.
$ws.Rows.Item(893).Insert()

# Column A holds text that looks like a date ("2026/02/27"); write it with a
# leading apostrophe so COM stores it as text instead of auto-converting to a
# date serial number, then reset the cell style so no stray number-format /
# quote-prefix style sticks around on the new cell.
$ws.Cells.Item(893, 1).Value = "'2026/02/27"
$ws.Cells.Item(893, 1).Style = "Normal"
$ws.Cells.Item(893, 2).Value = "金"
$ws.Cells.Item(893, 3).Value = 4
$ws.Cells.Item(893, 4).Value = 201
